$d = $word.ActiveDocument

# Anchor immediately before the documents final (trailing empty) paragraph,
# so new paragraphs get inserted ahead of it while it stays last & untouched.
$countBefore = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Last
$anchor = $d.Range($finalPara.Range.Start, $finalPara.Range.Start)

$numNew = 14
for ($i = 0; $i -lt $numNew; $i++) {
    $anchor.InsertParagraphBefore()
}

$fragments = @(
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br w:type="page"/></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="center"/><w:spacing w:before="200" w:after="100"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve">ANÁLISIS</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="both"/><w:spacing w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">PM10 (µg/m³)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="both"/><w:spacing w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">PM2.5 (µg/m³)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="both"/><w:spacing w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">SO₂ (PPB)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="both"/><w:spacing w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">NO₂ (PPB)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="both"/><w:spacing w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">O3 (PPB)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="both"/><w:spacing w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">CO (PPB)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="both"/><w:spacing w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">Velocidad del viento (m/s)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="both"/><w:spacing w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">Dirección del viento (º)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="both"/><w:spacing w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">Precipitación (mm)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="both"/><w:spacing w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">Presión atmosférica (mmHg)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="both"/><w:spacing w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">Temperatura (°C)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind/><w:jc w:val="both"/><w:spacing w:line="240"/></w:pPr><w:r><w:rPr><w:color w:val="002060"/><w:b/><w:bCs/><w:rFonts w:ascii="Microsoft YaHei UI" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Microsoft YaHei UI" w:cs="Microsoft YaHei UI"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">Humedad (%)</w:t></w:r></w:p>'
)

for ($i = 0; $i -lt $numNew; $i++) {
    $para = $d.Paragraphs.Item($countBefore + $i)
    $para.Range.InsertXML($fragments[$i])
}

Write-Output ("final paragraph count=" + $d.Paragraphs.Count)
